$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final state of the leaderboard (rows 2-23), reflecting updated
# total_registros counts and the resulting re-sort of tied/adjacent rows.
$data = @(
  @("FARFAN MONTOYA ROSSANA ISABEL", 62),
  @("DAVILA CORDOVA MARIBEL", 62),
  @("SANCHEZ ULLOA CESAR AUGUSTO", 57),
  @("RAMOS RAMOS HANDY JAIR", 57),
  @("FIORELA KEILY GUTIERREZ CRUZ", 55),
  @("CARBAJAL RAMOS JESUS MARINA", 54),
  @("GARCIA GUTIERREZ LUIS ARTURO", 54),
  @("ANGIE BELÉN RODRÍGUEZ ZAVALA", 52),
  @("GONZALES VALLE SEBASTIAN", 52),
  @("DELGADO DELGADO RONI", 51),
  @("OLIVA ALVA GOSSELYN NASSIRA", 51),
  @("BAZAN TEJADA JOSE VICENTE", 50),
  @("ARANEDA LOPEZ MARCO VIERI", 50),
  @("DANY DARWIN VILLACORTA SAAVEDRA", 49),
  @("VERDE LIZARRAGA DEYSI EUFEMIA", 48),
  @("CASTILLO QUEZADA DIEGO ALONSO", 46),
  @("GUERRA CALDERON ESTHEFANY NICOLLE", 45),
  @("JOSSY IVANA SUÁREZ ZAVALETA", 44),
  @("CARDENAS CAMPOJO MARY PAULA", 41),
  @("JAVE CHAVEZ ANGHELO MARTIN", 38),
  @("ALVITES CAMPOS SERGIO MARTIN", 20),
  @("ROCHA SIPIRAN JHORDAN ENRIQUE", 1)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row = $row + 1
}
